$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 28 — shifts the existing rows 28:63 down to 29:64
# (mirrors how a new daily price record gets prepended to this dataset).
$ws.Rows(28).Insert()

# Populate the newly inserted row 28 with the new record's data.
$ws.Cells.Item(28, 1).Value = 11
$ws.Cells.Item(28, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(28, 3).Value = "Bíobío"
$ws.Cells.Item(28, 4).Value = 44540
$ws.Cells.Item(28, 5).Value = 8
$ws.Cells.Item(28, 6).Value = "Fruta"
$ws.Cells.Item(28, 7).Value = 100101
$ws.Cells.Item(28, 8).Value = "Berries"
$ws.Cells.Item(28, 9).Value = 100101001
$ws.Cells.Item(28, 10).Value = "Arándano (blue)"
$ws.Cells.Item(28, 11).Value = "Sin especificar"
$ws.Cells.Item(28, 12).Value = "Primera"
$ws.Cells.Item(28, 13).Value = 250
$ws.Cells.Item(28, 14).Value = 3000
$ws.Cells.Item(28, 15).Value = 3500
$ws.Cells.Item(28, 16).Value = 3300
$ws.Cells.Item(28, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(28, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(28, 19).Value = 1650
$ws.Cells.Item(28, 20).Value = 2
